# Add a manual line break (as a separate run) at the end of each of the
# "Part 1: Photo effects" filter list items, right after the existing text
# run and before the paragraph mark.

$d = $word.ActiveDocument

$targets = @(
    "Exposure Filter",
    "Contrast Filter",
    "Saturation Filter",
    "Temperature Filter",
    "Solarization (White Clipping)"
)

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    foreach ($t in $targets) {
        if ($text -eq $t) {
            $r = $p.Range
            $r.Collapse(0)
            $r.InsertBreak(6)
        }
    }
}
